# "Fin du projet : restructuration"
# Append one more credential row (Ryan's own account) to the SHA1 table on
# Feuil1: row 21 = email / password / sha1-hash, with the email cell
# hyperlinked (mailto:) and styled like the rest of column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

$newRow = 21
$emailAddr = "ryan@ryanmalonzo.fr"

# Create the mailto: hyperlink on the new cell first (this also mints the
# worksheet relationship, rIdNN, that the hyperlink points at).
$ws.Hyperlinks.Add($ws.Cells.Item($newRow, 1), "mailto:$emailAddr")

# Fill in the new row's values.
$ws.Cells.Item($newRow, 1).Value = $emailAddr
$ws.Cells.Item($newRow, 2).Value = "Souris1234"
$ws.Cells.Item($newRow, 3).Value = "1f1d381923c82fc4e3df69c498053618eadf7a47"

# Match the hyperlink-cell styling used by every other row in column A
# (underlined "Lien hypertexte" look) by copying the format from the row
# above instead of relying on whatever default style Hyperlinks.Add applied.
$ws.Cells.Item($newRow - 1, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial($xlPasteFormats)

# Excel leaves the selection one cell below/right of the freshly entered
# data once you're done typing a row.
$ws.Range("C22").Select()
